$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47

function Set-TextCell($r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Text columns (Date, Time, Weekday, Week) - stored as literal text, not
# auto-converted to date/number serials, matching the existing rows.
Set-TextCell $row 1 "2023-06-15"
Set-TextCell $row 2 "15:07:00"
Set-TextCell $row 3 "Thursday"
Set-TextCell $row 4 "24"

# Numeric columns (city resale counts)
$ws.Cells.Item($row, 5).Value = 121554
$ws.Cells.Item($row, 6).Value = 132403
$ws.Cells.Item($row, 7).Value = 161738
$ws.Cells.Item($row, 8).Value = 132759
$ws.Cells.Item($row, 9).Value = 176424
$ws.Cells.Item($row, 10).Value = 114219
$ws.Cells.Item($row, 11).Value = 199204
$ws.Cells.Item($row, 12).Value = 223899
$ws.Cells.Item($row, 13).Value = 174319
$ws.Cells.Item($row, 14).Value = 102378
$ws.Cells.Item($row, 15).Value = 38987
$ws.Cells.Item($row, 16).Value = 34036
$ws.Cells.Item($row, 17).Value = 51655
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36447
$ws.Cells.Item($row, 20).Value = -1
